$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).ClearFormats()
}

Set-TextValue "D2" "29.135.68"
$ws.Range("E2").Value = "  -1.85%  "
Set-TextValue "D3" "1.838.68"
$ws.Range("E3").Value = "  -1.41%  "
Set-TextValue "D4" "0.9994"
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue "D5" "239.84"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("E7").Value = "  -0.07%  "
Set-TextValue "D8" "0.2990"
$ws.Range("E8").Value = "  -2.96%  "
Set-TextValue "D9" "0.07451"
$ws.Range("E9").Value = "  -4.16%  "
Set-TextValue "D10" "23.22"
$ws.Range("E10").Value = "  -2.27%  "
Set-TextValue "D11" "0.07645"
$ws.Range("E11").Value = "  -2.38%  "
Set-TextValue "D12" "1.842.39"
$ws.Range("E12").Value = "  -1.20%  "
Set-TextValue "D13" "5.031"
$ws.Range("E13").Value = "  -2.85%  "
Set-TextValue "D14" "0.6807"
$ws.Range("E14").Value = "  -2.23%  "
Set-TextValue "D15" "87.43"
$ws.Range("E15").Value = "  -5.93%  "
Set-TextValue "D16" "6.153"
$ws.Range("E16").Value = "  -7.29%  "
Set-TextValue "D17" "29.145.44"
$ws.Range("E17").Value = "  -1.82%  "
Set-TextValue "D18" "0.000008212"
$ws.Range("E18").Value = "  -2.18%  "
Set-TextValue "D19" "2.092.09"
$ws.Range("E19").Value = "  -1.28%  "
Set-TextValue "D20" "230.50"
$ws.Range("E20").Value = "  -5.43%  "
Set-TextValue "D21" "12.50"
$ws.Range("E21").Value = "  -2.50%  "
Set-TextValue "D22" "0.9998"
$ws.Range("E22").Value = "  -0.05%  "
Set-TextValue "D23" "7.346"
$ws.Range("E23").Value = "  -4.09%  "
Set-TextValue "D24" "1.0000"
$ws.Range("E24").Value = "  -0.06%  "
Set-TextValue "D25" "161.15"
$ws.Range("E25").Value = "  +0.59%  "
Set-TextValue "D26" "0.1428"
$ws.Range("E26").Value = "  -6.26%  "
Set-TextValue "D27" "8.706"
$ws.Range("E27").Value = "  -2.96%  "
$ws.Range("E28").Value = "  -1.78%  "
Set-TextValue "D29" "1.503"
$ws.Range("E29").Value = "  -2.82%  "
Set-TextValue "D30" "4.252"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("E32").Value = "  -0.50%  "
Set-TextValue "D33" "0.05343"
$ws.Range("E33").Value = "  +4.51%  "
Set-TextValue "D34" "0.7539"
$ws.Range("E34").Value = "  -4.47%  "
Set-TextValue "D35" "1.847"
$ws.Range("E35").Value = "  -3.65%  "
Set-TextValue "D36" "1.133"
$ws.Range("E36").Value = "  -2.44%  "
Set-TextValue "D37" "2.683"
$ws.Range("E37").Value = "  -0.33%  "
Set-TextValue "D38" "1.311.79"
$ws.Range("E38").Value = "  -2.28%  "
Set-TextValue "D39" "0.01826"
$ws.Range("E39").Value = "  -3.15%  "
Set-TextValue "D40" "2.716"
$ws.Range("E40").Value = "  -0.98%  "
Set-TextValue "D41" "0.9435"
$ws.Range("E41").Value = "  -1.87%  "
Set-TextValue "D42" "6.066"
$ws.Range("E42").Value = "  +0.40%  "
Set-TextValue "D43" "104.83"
$ws.Range("E43").Value = "  -1.78%  "
Set-TextValue "D44" "0.9991"
Set-TextValue "D45" "0.07961"
$ws.Range("E45").Value = "  +25.27%  "
Set-TextValue "D46" "1.987.57"
$ws.Range("E46").Value = "  -1.32%  "
Set-TextValue "D47" "0.5180"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("E48").Value = "  -4.08%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "1.773"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D50" "64.00"
$ws.Range("E50").Value = "  -1.96%  "
Set-TextValue "D51" "9.429"
$ws.Range("E51").Value = "  -3.81%  "
